$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Row 6 is the "Day 2" row. Rename the placeholder title and fill in the
# actual times/rank now that Day 2 ("1202 Program Alarm") is done.
$ws.Range("B6").Value = "Day 2: 1202 Program Alarm"
$ws.Range("C6").Value = 0.006782407407407408
$ws.Range("E6").Value = 0.012627314814814815
$ws.Range("F6").Value = 0.009444444444444445
$ws.Range("H6").Value = "3rd"

# Move the active selection down to H7, matching where the user ended up.
$ws.Range("H7").Select()

$excel.Calculate()
